$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 5263855
$ws.Range("I33").Value = 10000545
$ws.Range("J33").Value = 866.7778
$ws.Range("K33").Value = 10000545
$ws.Range("L33").Value = 866.7778
$ws.Range("M33").Value = -10000316
$ws.Range("N33").Value = -1324.7778

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 1772.3334
$ws.Range("J46").Value = 3000
$ws.Range("L46").Value = 9000
$ws.Range("N46").Value = -9238

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H60").Value = 1772.3334
$ws.Range("J60").Value = 3000
$ws.Range("L60").Value = 9000
$ws.Range("N60").Value = -9968

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4332.3335
$ws.Range("I76").Value = 2800
$ws.Range("K76").Value = 2800
$ws.Range("M76").Value = -2485

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 4332.3335
$ws.Range("I79").Value = 2800
$ws.Range("K79").Value = 2800
$ws.Range("M79").Value = -1708

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 6976
$ws.Range("I98").Value = 6976
$ws.Range("K98").Value = 6976
$ws.Range("M98").Value = -5478

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 3537.7058
$ws.Range("I106").Value = 3484.6428
$ws.Range("K106").Value = 3484.6428
$ws.Range("M106").Value = -2853.6428

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 2788.88
$ws.Range("J107").Value = 4731.3335
$ws.Range("L107").Value = 4731.3335
$ws.Range("N107").Value = -8571.333500000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 64697.4
$ws.Range("I113").Value = 83139.21000000001
$ws.Range("K113").Value = 83139.21000000001
$ws.Range("M113").Value = -79885.21000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 9562.941000000001
$ws.Range("I116").Value = 9820.571
$ws.Range("J116").Value = 9382.6
$ws.Range("K116").Value = 9820.571
$ws.Range("L116").Value = 9382.6
$ws.Range("M116").Value = -6378.571
$ws.Range("N116").Value = -16266.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 6976
$ws.Range("I122").Value = 6976
$ws.Range("K122").Value = 20928
$ws.Range("M122").Value = -18478

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 7622529.5
$ws.Range("I132").Value = 9037257
$ws.Range("J132").Value = 144683.86
$ws.Range("K132").Value = 27111771
$ws.Range("L132").Value = 434051.58
$ws.Range("M132").Value = -27109241
$ws.Range("N132").Value = -439111.58

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3799.2856
$ws.Range("I63").Value = 3795.3333
$ws.Range("K63").Value = 3795.3333
$ws.Range("M63").Value = -3109.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3799.2856
$ws.Range("I66").Value = 3795.3333
$ws.Range("K66").Value = 18976.6665
$ws.Range("M66").Value = -15544.6665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1598.8889
$ws.Range("I97").Value = 1633.125
$ws.Range("K97").Value = 1633.125
$ws.Range("M97").Value = -1137.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 50786.57
$ws.Range("I122").Value = 61694.41
$ws.Range("K122").Value = 185083.23
$ws.Range("M122").Value = -182633.23

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2512.4324
$ws.Range("I132").Value = 2171.5151
$ws.Range("J132").Value = 5325
$ws.Range("K132").Value = 6514.5453
$ws.Range("L132").Value = 15975
$ws.Range("M132").Value = -3984.5453
$ws.Range("N132").Value = -21035

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1736.8889
$ws.Range("I94").Value = 1946
$ws.Range("J94").Value = 1005
$ws.Range("K94").Value = 1946
$ws.Range("L94").Value = 1005
$ws.Range("M94").Value = -1495
$ws.Range("N94").Value = -1907

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3679
$ws.Range("J107").Value = 1531.6666
$ws.Range("L107").Value = 1531.6666
$ws.Range("N107").Value = -5371.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1366.5778
$ws.Range("I134").Value = 985.9211
$ws.Range("J134").Value = 3433
$ws.Range("K134").Value = 2957.7633
$ws.Range("L134").Value = 10299
$ws.Range("M134").Value = -422.7633000000001
$ws.Range("N134").Value = -15369

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 18779262
$ws.Range("I31").Value = 2436071.5
$ws.Range("J31").Value = 125010010
$ws.Range("K31").Value = 2436071.5
$ws.Range("L31").Value = 125010010
$ws.Range("M31").Value = -2435776.5
$ws.Range("N31").Value = -125010600

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 18779262
$ws.Range("I34").Value = 2436071.5
$ws.Range("J34").Value = 125010010
$ws.Range("K34").Value = 2436071.5
$ws.Range("L34").Value = 125010010
$ws.Range("M34").Value = -2435869.5
$ws.Range("N34").Value = -125010414

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H95").Value = 127749.4
$ws.Range("J95").Value = 154187
$ws.Range("L95").Value = 154187
$ws.Range("N95").Value = -159679

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H96").Value = 7890.8
$ws.Range("J96").Value = 7890.8
$ws.Range("L96").Value = 7890.8
$ws.Range("N96").Value = -13382.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 12495.897
$ws.Range("I134").Value = 12275.17
$ws.Range("K134").Value = 36825.51
$ws.Range("M134").Value = -34290.51

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 2255.5715
$ws.Range("I50").Value = 72.25
$ws.Range("K50").Value = 216.75
$ws.Range("M50").Value = 264.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H53").Value = 2255.5715
$ws.Range("I53").Value = 72.25
$ws.Range("K53").Value = 216.75
$ws.Range("M53").Value = 264.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2632805.5
$ws.Range("J68").Value = 4168117.2
$ws.Range("L68").Value = 12504351.6
$ws.Range("N68").Value = -12505973.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 2632805.5
$ws.Range("J71").Value = 4168117.2
$ws.Range("L71").Value = 37513054.8
$ws.Range("N71").Value = -37521166.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 25641850
$ws.Range("I107").Value = 33333796
$ws.Range("J107").Value = 2028.3334
$ws.Range("K107").Value = 100001388
$ws.Range("L107").Value = 6085.0002
$ws.Range("M107").Value = -99999468
$ws.Range("N107").Value = -9925.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 418
$ws.Range("I113").Value = 496.66666
$ws.Range("K113").Value = 1489.99998
$ws.Range("M113").Value = 680.0000199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 2453.4
$ws.Range("J122").Value = 3401
$ws.Range("L122").Value = 30609
$ws.Range("N122").Value = -35509

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 17975.031
$ws.Range("J131").Value = 2675.827
$ws.Range("L131").Value = 8027.481000000001
$ws.Range("N131").Value = -18107.481

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2150
$ws.Range("I102").Value = 1825.0625
$ws.Range("K102").Value = 1825.0625
$ws.Range("M102").Value = -203.0625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H135").Value = 124500
$ws.Range("J135").Value = 124500
$ws.Range("L135").Value = 124500
$ws.Range("N135").Value = -134640

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2681.5454
$ws.Range("I126").Value = 2499.6667
$ws.Range("K126").Value = 7499.000100000001
$ws.Range("M126").Value = -5029.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3458.1724
$ws.Range("I132").Value = 3545.9644
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 10637.8932
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -8107.893199999999
$ws.Range("N132").Value = -8060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1582.92
$ws.Range("I136").Value = 1265.1578
$ws.Range("J136").Value = 2589.1667
$ws.Range("K136").Value = 3795.4734
$ws.Range("L136").Value = 7767.500100000001
$ws.Range("M136").Value = -1245.4734
$ws.Range("N136").Value = -12867.5001

